$wb = $excel.ActiveWorkbook

# --- Layer0 sheet ---
$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.3215772883448896
$ws0.Range("C2").Value = -0.8878702236976402
$ws0.Range("B3").Value = -1.095000745142599
$ws0.Range("C3").Value = 0.5409518889988462
$ws0.Range("B4").Value = -1.079928076358642
$ws0.Range("C4").Value = 1.245558183491608

# --- Layer1 sheet ---
$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -1.405220245303278
$ws1.Range("C2").Value = 0.2984837608913679
$ws1.Range("B3").Value = -1.101721460397068
$ws1.Range("C3").Value = -0.7389973702544742
$ws1.Range("B4").Value = 1.285324352609827
$ws1.Range("C4").Value = -0.3118840604879607
